# Auto-generated Excel COM-interop script to apply scheduled-runner data refresh
# to the Unicorn_Profits workbook (updates currentAveragePrice / Leve price / profit
# columns H-N across 8 sheets, 36 rows total, based on upstream market data).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2018.5942
$ws.Range("I40").Value = 1115.8334
$ws.Range("K40").Value = 1115.8334
$ws.Range("M40").Value = -940.8334
# Row 64
$ws.Range("H64").Value = 3851.2083
$ws.Range("I64").Value = 3566.0833
$ws.Range("J64").Value = 4706.5835
$ws.Range("K64").Value = 3566.0833
$ws.Range("L64").Value = 4706.5835
$ws.Range("M64").Value = -3318.0833
$ws.Range("N64").Value = -5202.5835
# Row 67
$ws.Range("H67").Value = 3851.2083
$ws.Range("I67").Value = 3566.0833
$ws.Range("J67").Value = 4706.5835
$ws.Range("K67").Value = 3566.0833
$ws.Range("L67").Value = 4706.5835
$ws.Range("M67").Value = -2708.0833
$ws.Range("N67").Value = -6422.5835
# Row 74
$ws.Range("H74").Value = 5058.8237
$ws.Range("I74").Value = 4912
$ws.Range("J74").Value = 5268.5713
$ws.Range("K74").Value = 4912
$ws.Range("L74").Value = 5268.5713
$ws.Range("M74").Value = -3976
$ws.Range("N74").Value = -7140.5713
# Row 76
$ws.Range("H76").Value = 14712422
$ws.Range("I76").Value = 20840538
$ws.Range("J76").Value = 4946.4
$ws.Range("K76").Value = 20840538
$ws.Range("L76").Value = 4946.4
$ws.Range("M76").Value = -20840223
$ws.Range("N76").Value = -5576.4
# Row 77
$ws.Range("H77").Value = 5058.8237
$ws.Range("I77").Value = 4912
$ws.Range("J77").Value = 5268.5713
$ws.Range("K77").Value = 24560
$ws.Range("L77").Value = 26342.8565
$ws.Range("M77").Value = -19880
$ws.Range("N77").Value = -35702.85649999999
# Row 79
$ws.Range("H79").Value = 14712422
$ws.Range("I79").Value = 20840538
$ws.Range("J79").Value = 4946.4
$ws.Range("K79").Value = 20840538
$ws.Range("L79").Value = 4946.4
$ws.Range("M79").Value = -20839446
$ws.Range("N79").Value = -7130.4
# Row 116
$ws.Range("H116").Value = 55839.35
$ws.Range("I116").Value = 73287.13
$ws.Range("K116").Value = 73287.13
$ws.Range("M116").Value = -69845.13
# Row 132
$ws.Range("H132").Value = 3688.2327
$ws.Range("I132").Value = 1727.6111
$ws.Range("J132").Value = 13771.429
$ws.Range("K132").Value = 5182.8333
$ws.Range("L132").Value = 41314.287
$ws.Range("M132").Value = -2652.8333
$ws.Range("N132").Value = -46374.287

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2230
$ws.Range("I2").Value = 3445
$ws.Range("J2").Value = 1258
$ws.Range("K2").Value = 3445
$ws.Range("L2").Value = 1258
$ws.Range("M2").Value = -3332
$ws.Range("N2").Value = -1484
# Row 6
$ws.Range("H6").Value = 2000000
$ws.Range("I6").Value = 2000000
$ws.Range("K6").Value = 2000000
$ws.Range("M6").Value = -1999827
# Row 32
$ws.Range("H32").Value = 1285526.2
$ws.Range("I32").Value = 1491724.9
$ws.Range("J32").Value = 2512.5557
$ws.Range("K32").Value = 1491724.9
$ws.Range("L32").Value = 2512.5557
$ws.Range("M32").Value = -1491437.9
$ws.Range("N32").Value = -3086.5557
# Row 63
$ws.Range("H63").Value = 2646
$ws.Range("I63").Value = 2689.0667
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 2689.0667
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -2003.0667
$ws.Range("N63").Value = -3372
# Row 66
$ws.Range("H66").Value = 2646
$ws.Range("I66").Value = 2689.0667
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 13445.3335
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -10013.3335
$ws.Range("N66").Value = -16864
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
# Row 116
$ws.Range("H116").Value = 2230
$ws.Range("I116").Value = 3445
$ws.Range("J116").Value = 1258
$ws.Range("K116").Value = 3445
$ws.Range("L116").Value = 1258
$ws.Range("M116").Value = -1151
$ws.Range("N116").Value = -5846

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2230
$ws.Range("I3").Value = 3445
$ws.Range("J3").Value = 1258
$ws.Range("K3").Value = 3445
$ws.Range("L3").Value = 1258
$ws.Range("M3").Value = -3331
$ws.Range("N3").Value = -1486
# Row 22
$ws.Range("H22").Value = 337.875
$ws.Range("I22").Value = 302.2
$ws.Range("J22").Value = 397.33334
$ws.Range("K22").Value = 302.2
$ws.Range("L22").Value = 397.33334
$ws.Range("M22").Value = -129.2
$ws.Range("N22").Value = -743.33334
# Row 40
$ws.Range("H40").Value = 30000
$ws.Range("J40").Value = 30000
$ws.Range("L40").Value = 30000
$ws.Range("N40").Value = -30530
# Row 86
$ws.Range("H86").Value = 3361.9333
$ws.Range("I86").Value = 4672.7334
$ws.Range("J86").Value = 2051.1333
$ws.Range("K86").Value = 4672.7334
$ws.Range("L86").Value = 2051.1333
$ws.Range("M86").Value = -3549.7334
$ws.Range("N86").Value = -4297.1333
# Row 89
$ws.Range("H89").Value = 3361.9333
$ws.Range("I89").Value = 4672.7334
$ws.Range("J89").Value = 2051.1333
$ws.Range("K89").Value = 23363.667
$ws.Range("L89").Value = 10255.6665
$ws.Range("M89").Value = -17747.667
$ws.Range("N89").Value = -21487.6665
# Row 99
$ws.Range("H99").Value = 12715.714
$ws.Range("I99").Value = 11402
$ws.Range("J99").Value = 16000
$ws.Range("K99").Value = 11402
$ws.Range("L99").Value = 16000
$ws.Range("M99").Value = -9904
$ws.Range("N99").Value = -18996

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 74
$ws.Range("H74").Value = 12589
$ws.Range("J74").Value = 12589
$ws.Range("L74").Value = 12589
$ws.Range("N74").Value = -14337
# Row 77
$ws.Range("H77").Value = 12589
$ws.Range("J77").Value = 12589
$ws.Range("L77").Value = 37767
$ws.Range("N77").Value = -46503

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1112.3334
$ws.Range("I5").Value = 403.6
$ws.Range("J5").Value = 1998.25
$ws.Range("K5").Value = 1210.8
$ws.Range("L5").Value = 5994.75
$ws.Range("M5").Value = -1098.8
$ws.Range("N5").Value = -6218.75
# Row 38
$ws.Range("H38").Value = 72.7
$ws.Range("I38").Value = 48.46154
$ws.Range("J38").Value = 117.71429
$ws.Range("K38").Value = 145.38462
$ws.Range("L38").Value = 353.14287
$ws.Range("M38").Value = 201.61538
$ws.Range("N38").Value = -1047.14287
# Row 98
$ws.Range("H98").Value = 5889.5557
$ws.Range("I98").Value = 217.66667
$ws.Range("J98").Value = 17233.334
$ws.Range("K98").Value = 653.00001
$ws.Range("L98").Value = 51700.00199999999
$ws.Range("M98").Value = 844.99999
$ws.Range("N98").Value = -54696.00199999999
# Row 135
$ws.Range("H135").Value = 1112.3334
$ws.Range("I135").Value = 403.6
$ws.Range("J135").Value = 1998.25
$ws.Range("K135").Value = 3632.4
$ws.Range("L135").Value = 17984.25
$ws.Range("M135").Value = -1097.4
$ws.Range("N135").Value = -23054.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 387.5
$ws.Range("I46").Value = 387.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 387.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -231.5
$ws.Range("N46").ClearContents()
# Row 57
$ws.Range("H57").Value = 10955
$ws.Range("I57").Value = 10955
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 10955
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -10135
$ws.Range("N57").ClearContents()
# Row 80
$ws.Range("H80").Value = 4040.34
$ws.Range("I80").Value = 4413.237
$ws.Range("J80").Value = 2859.5
$ws.Range("K80").Value = 4413.237
$ws.Range("L80").Value = 2859.5
$ws.Range("M80").Value = -3415.237
$ws.Range("N80").Value = -4855.5
# Row 83
$ws.Range("H83").Value = 4040.34
$ws.Range("I83").Value = 4413.237
$ws.Range("J83").Value = 2859.5
$ws.Range("K83").Value = 22066.185
$ws.Range("L83").Value = 14297.5
$ws.Range("M83").Value = -17074.185
$ws.Range("N83").Value = -24281.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 105
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988
# Row 122
$ws.Range("H122").Value = 2599.2173
$ws.Range("I122").Value = 2286.1333
$ws.Range("K122").Value = 6858.3999
$ws.Range("M122").Value = -4408.3999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

